$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string (rich text) cell updates ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Plain numeric / simple value updates ---
$ws.Range("F14").Value = 1
$ws.Range("M14").Value = 300
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 23
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = 4.545454545454
$ws.Range("L15").Value = -4.166666666666
$ws.Range("M15").Value = 64.285714285714
$ws.Range("N15").Value = -45.238095238095
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -6.666666666666
$ws.Range("I16").Value = 270
$ws.Range("J16").Value = 306
$ws.Range("K16").Value = -11.764705882352
$ws.Range("L16").Value = 27.962085308056
$ws.Range("M16").Value = 8.433734939759
$ws.Range("N16").Value = -64
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 41.666666666666
$ws.Range("F17").Value = 65
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 54.761904761904
$ws.Range("I17").Value = 469
$ws.Range("J17").Value = 429
$ws.Range("K17").Value = 9.324009324009
$ws.Range("L17").Value = 49.36305732484
$ws.Range("M17").Value = 83.203125
$ws.Range("N17").Value = -3.497942386831
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -42.857142857142
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 41.176470588235
$ws.Range("I18").Value = 164
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = -5.747126436781
$ws.Range("L18").Value = 19.70802919708
$ws.Range("M18").Value = 0.613496932515
$ws.Range("N18").Value = -81.675977653631
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 36.363636363636
$ws.Range("I19").Value = 338
$ws.Range("J19").Value = 307
$ws.Range("K19").Value = 10.097719869706
$ws.Range("L19").Value = 11.551155115511
$ws.Range("M19").Value = 52.252252252252
$ws.Range("N19").Value = 8.333333333333
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 233
$ws.Range("J20").Value = 196
$ws.Range("K20").Value = 18.877551020408
$ws.Range("L20").Value = 94.166666666666
$ws.Range("M20").Value = 187.654320987654
$ws.Range("N20").Value = -37.700534759358
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 190
$ws.Range("G21").Value = 154
$ws.Range("H21").Value = 23.376623376623
$ws.Range("I21").Value = 1509
$ws.Range("J21").Value = 1446
$ws.Range("K21").Value = 4.356846473029
$ws.Range("L21").Value = 34.973166368515
$ws.Range("M21").Value = 52.732793522267
$ws.Range("N21").Value = -47.676837725381
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = -19.047619047619
$ws.Range("L23").Value = 41.666666666666
$ws.Range("M23").Value = 54.545454545454
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -62.5
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = -39.805825242718
$ws.Range("I24").Value = 690
$ws.Range("J24").Value = 841
$ws.Range("K24").Value = -17.9548156956
$ws.Range("L24").Value = 20.840630472854
$ws.Range("M24").Value = 1.024890190336
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 41.176470588235
$ws.Range("F25").Value = 85
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = 39.344262295082
$ws.Range("I25").Value = 615
$ws.Range("J25").Value = 610
$ws.Range("K25").Value = 0.819672131147
$ws.Range("L25").Value = 24.746450304259
$ws.Range("M25").Value = -8.888888888888
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = -57.142857142857
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 43
$ws.Range("K26").Value = -16.279069767441
$ws.Range("L26").Value = 0
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 52
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = -8.771929824561
$ws.Range("L27").Value = 44.444444444444
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("M28").Value = -16.129032258064
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("M29").Value = -26.923076923076

# --- Cells whose type changes (number <-> text) need a style fix-up ---
# Strategy: set the value, then copy number-format/style from a donor cell
# that already carries the desired style, via PasteSpecial(xlPasteFormats).
$ws.Range("G22").Value = "'0"
$ws.Range("F22").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null
$ws.Range("H22").Value = "***.*"
$ws.Range("E22").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 1
$ws.Range("C26").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 1
$ws.Range("C26").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = 2
$ws.Range("D26").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = -50
$ws.Range("E26").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = "'0"
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = "'0"
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "***.*"
$ws.Range("E30").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Value = "'0"
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Value = "'0"
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = "***.*"
$ws.Range("E30").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
